$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cargos_tarefas")

$ws.Range('C16').Value = 'Superior de Macau'
$ws.Range('D16').Value = 'jesuita-cargo'
$ws.Range('C17').Value = 'Funda a residência jesuíta de Macau'
$ws.Range('D17').Value = 'jesuita-tarefa'
$ws.Range('A46').Value = 'deh-nicolau-pimenta'
$ws.Range('B46').Value = 'Nicolau Pimenta'
$ws.Range('C46').Value = 'Visitador das Índias Orientais'
$ws.Range('D46').Value = 'jesuita-cargo'
$ws.Range('A47').Value = 'deh-joao-rodrigues-girao'
$ws.Range('B47').Value = 'João Rodrigues Girão'
$ws.Range('C47').Value = 'Carta anual do Japão'
$ws.Range('D47').Value = 'jesuita-tarefa'
$ws.Range('A60').Value = 'deh-mateus-de-couros'
$ws.Range('B60').Value = 'Mateus de Couros'
$ws.Range('C60').Value = 'Provincial do Japão e da China'
$ws.Range('A61').Value = 'deh-francisco-pacheco'
$ws.Range('B61').Value = 'Francisco Pacheco'
$ws.Range('C61').Value = 'Provincial do Japão'
$ws.Range('A63').Value = 'deh-manuel-dias-o-novo'
$ws.Range('B63').Value = 'Manuel Dias, o Novo'
$ws.Range('C63').Value = 'Vice-provincial da China'
$ws.Range('A64').Value = 'deh-francisco-pacheco'
$ws.Range('B64').Value = 'Francisco Pacheco'
$ws.Range('C64').Value = 'Superior da missão da China'
$ws.Range('D64').Value = 'jesuita-cargo'
$ws.Range('A65').Value = 'deh-gaspar-do-amaral'
$ws.Range('B65').Value = 'Gaspar do Amaral'
$ws.Range('C65').Value = 'Faz mais de 40000 baptismos no Tonquim'
$ws.Range('D65').Value = 'jesuita-tarefa'
$ws.Range('A71').Value = 'deh-mateus-de-couros'
$ws.Range('B71').Value = 'Mateus de Couros'
$ws.Range('C71').Value = 'Provincial do Japão e da China'
$ws.Range('C72').Value = 'Visitador das províncias de Goa e do Malabar'
$ws.Range('A73').Value = 'deh-andre-palmeiro'
$ws.Range('B73').Value = 'André Palmeiro'
$ws.Range('C73').Value = 'Visitador do Japão e da China'
$ws.Range('D73').Value = 'jesuita-cargo'
$ws.Range('A74').Value = 'deh-joao-rodrigues-girao'
$ws.Range('B74').Value = 'João Rodrigues Girão'
$ws.Range('C74').Value = 'Carta anual do Japão'
$ws.Range('D74').Value = 'jesuita-tarefa'
$ws.Range('A76').Value = 'deh-gaspar-do-amaral'
$ws.Range('B76').Value = 'Gaspar do Amaral'
$ws.Range('C76').Value = 'Superior da missão do Tonquim'
$ws.Range('A77').Value = 'deh-antonio-de-andrade'
$ws.Range('B77').Value = 'António de Andrade'
$ws.Range('C77').Value = 'Provincial de Goa'
$ws.Range('A81').Value = 'deh-manuel-dias-o-novo'
$ws.Range('B81').Value = 'Manuel Dias, o Novo'
$ws.Range('C81').Value = 'Vice-provincial da China'
$ws.Range('A83').Value = 'deh-andre-palmeiro'
$ws.Range('B83').Value = 'André Palmeiro'
$ws.Range('C83').Value = 'Visitador do Japão e da China'
$ws.Range('A95').Value = 'deh-manuel-dias-o-novo'
$ws.Range('B95').Value = 'Manuel Dias, o Novo'
$ws.Range('C95').Value = 'Vice-provincial da China'
$ws.Range('A96').Value = 'deh-francisco-furtado'
$ws.Range('B96').Value = 'Francisco Furtado'
$ws.Range('C96').Value = 'Superior das seis residências do Norte'
$ws.Range('A105').Value = 'deh-simao-da-cunha'
$ws.Range('B105').Value = 'Simão da Cunha'
$ws.Range('C105').Value = 'Procurador da Vice-província da China'
$ws.Range('A106').Value = 'deh-matias-da-maia'
$ws.Range('B106').Value = 'Matias da Maia'
$ws.Range('C106').Value = 'Provincial do Japão'
$ws.Range('A107').Value = 'deh-sebastiao-da-maia'
$ws.Range('B107').Value = 'Sebastião da Maia'
$ws.Range('C107').Value = 'Visitador do Japão e da China'
$ws.Range('A123').Value = 'deh-miguel-do-amaral'
$ws.Range('B123').Value = 'Miguel do Amaral'
$ws.Range('C123').Value = 'Procurador da Vice-província da China'
$ws.Range('D123').Value = 'jesuita-cargo'
$ws.Range('A124').Value = 'deh-manuel-de-sa'
$ws.Range('B124').Value = 'Manuel de Sá'
$ws.Range('C124').Value = 'Enviado ao rei'
$ws.Range('D124').Value = 'jesuita-tarefa'
$ws.Range('A128').Value = 'deh-leonardo-teixeira'
$ws.Range('B128').Value = 'Leonardo Teixeira'
$ws.Range('C128').Value = 'Superior'
$ws.Range('A129').Value = 'deh-antonio-da-silva'
$ws.Range('B129').Value = 'António da Silva'
$ws.Range('C129').Value = 'Reitor de Nanquim'
$ws.Range('A132').Value = 'deh-manuel-osorio-i'
$ws.Range('B132').Value = 'Manuel Osório'
$ws.Range('C132').Value = 'Reitor do colégio de Macau'
$ws.Range('D132').Value = 'jesuita-cargo'
$ws.Range('C133').Value = 'Vice-provincial do Japão'
$ws.Range('A134').Value = 'deh-miguel-do-amaral'
$ws.Range('B134').Value = 'Miguel do Amaral'
$ws.Range('C134').Value = 'Provincial do Japão'
$ws.Range('A135').Value = 'deh-leonardo-teixeira'
$ws.Range('B135').Value = 'Leonardo Teixeira'
$ws.Range('C135').Value = 'Sai da Companhia'
$ws.Range('D135').Value = 'jesuita-tarefa'
$ws.Range('A137').Value = 'deh-manuel-osorio-i'
$ws.Range('B137').Value = 'Manuel Osório'
$ws.Range('C137').Value = 'Vice-provincial do Japão'
$ws.Range('C138').Value = 'Procurador da Província do Japão em Macau'
$ws.Range('A139').Value = 'deh-francisco-pinto-i'
$ws.Range('B139').Value = 'Francisco Pinto'
$ws.Range('C139').Value = 'Provincial do Japão'
$ws.Range('A143').Value = 'deh-miguel-do-amaral'
$ws.Range('B143').Value = 'Miguel do Amaral'
$ws.Range('A144').Value = 'deh-francisco-pinto-i'
$ws.Range('B144').Value = 'Francisco Pinto'
$ws.Range('A146').Value = 'deh-manuel-mendes'
$ws.Range('B146').Value = 'Manuel Mendes'
$ws.Range('A147').Value = 'deh-antonio-da-silva'
$ws.Range('B147').Value = 'António da Silva'
$ws.Range('A151').Value = 'deh-miguel-do-amaral'
$ws.Range('B151').Value = 'Miguel do Amaral'
$ws.Range('C151').Value = 'Provincial do Japão'
$ws.Range('C152').Value = 'Superior da missão do Tonquim'
$ws.Range('A153').Value = 'deh-estanislau-machado'
$ws.Range('B153').Value = 'Estanislau Machado'
$ws.Range('A158').Value = 'deh-miguel-do-amaral'
$ws.Range('B158').Value = 'Miguel do Amaral'
$ws.Range('C158').Value = 'Procurador da Vice-província da China'
$ws.Range('A159').Value = 'deh-manuel-mendes'
$ws.Range('B159').Value = 'Manuel Mendes'
$ws.Range('C159').Value = 'Vice-provincial da China'
$ws.Range('A160').Value = 'deh-francisco-pinto-i'
$ws.Range('B160').Value = 'Francisco Pinto'
$ws.Range('C160').Value = 'Vice-provincial do Japão'
$ws.Range('A161').Value = 'deh-francisco-pinto-i'
$ws.Range('B161').Value = 'Francisco Pinto'
$ws.Range('C161').Value = 'Vice-visitador do Japão'
